$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-12 (columns A-G)
$arr = New-Object 'object[,]' 11,7

# Row 2 (A2:G2)
$arr[0,0] = 1
$arr[0,1] = 221
$arr[0,2] = "Purse Gucci XL "
$arr[0,3] = "Branded Purse XL"
$arr[0,4] = "Fashion"
$arr[0,5] = 280.5
$arr[0,6] = 15

# Row 3 (A3:G3)
$arr[1,0] = 2
$arr[1,1] = 111
$arr[1,2] = "Laptopn I7 "
$arr[1,3] = "High Spec Laptop"
$arr[1,4] = "Electronics"
$arr[1,5] = 990.9
$arr[1,6] = 7

# Row 4 (A4:G4)
$arr[2,0] = 3
$arr[2,1] = 112
$arr[2,2] = "Laptopn I3 "
$arr[2,3] = "Mid Spec Laptop"
$arr[2,4] = "Electronics"
$arr[2,5] = 590.9
$arr[2,6] = 22

# Row 5 (A5:G5)
$arr[3,0] = 4
$arr[3,1] = 114
$arr[3,2] = "Laptopn I9 "
$arr[3,3] = "Local"
$arr[3,4] = "Electronics"
$arr[3,5] = 999.9
$arr[3,6] = 1

# Row 6 (A6:G6)
$arr[4,0] = 5
$arr[4,1] = 224
$arr[4,2] = "Shoes Addias "
$arr[4,3] = "Premium Shoes XL"
$arr[4,4] = "Shoes"
$arr[4,5] = 999.9
$arr[4,6] = 3

# Row 7 (A7:G7)
$arr[5,0] = 6
$arr[5,1] = 1
$arr[5,2] = "Item 1"
$arr[5,3] = "Description for Item 1"
$arr[5,4] = "Fashion"
$arr[5,5] = 10
$arr[5,6] = 5

# Row 8 (A8:G8)
$arr[6,0] = 7
$arr[6,1] = 3
$arr[6,2] = "Item 3"
$arr[6,3] = "Description for Item 3"
$arr[6,4] = "Fashion"
$arr[6,5] = 10
$arr[6,6] = 50

# Row 9 (A9:G9)
$arr[7,0] = 8
$arr[7,1] = 4
$arr[7,2] = "Item 4"
$arr[7,3] = "Description for Item 4"
$arr[7,4] = "Fashion"
$arr[7,5] = 10
$arr[7,6] = 50

# Row 10 (A10:G10)
$arr[8,0] = 9
$arr[8,1] = 10
$arr[8,2] = "Laptop"
$arr[8,3] = "A high-performance laptop"
$arr[8,4] = "Electronics"
$arr[8,5] = 999.99
$arr[8,6] = 10

# Row 11 (A11:G11)
$arr[9,0] = 10
$arr[9,1] = 11
$arr[9,2] = "Smartphone"
$arr[9,3] = "A latest model smartphone"
$arr[9,4] = "Electronics"
$arr[9,5] = 699.99
$arr[9,6] = 25

# Row 12 (A12:G12)
$arr[10,0] = 11
$arr[10,1] = 12
$arr[10,2] = "Desk Chair"
$arr[10,3] = "An ergonomic desk chair"
$arr[10,4] = "Furniture"
$arr[10,5] = 149.99
$arr[10,6] = 15

$ws.Range("A2:G12").Value = $arr
